$d = $word.ActiveDocument

# The last paragraph in the document body currently ends with the
# "-CL - Data Structres/Input Buffers..." bullet. Add a new bullet
# paragraph right after it containing the new idea text, matching the
# formatting of the preceding paragraphs (inherited automatically by
# InsertParagraphAfter from the paragraph it splits off of).
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.InsertAfter("- Platform/API agnostic back-end, to open up support for other graphic apis and platforms, including uwp/xbox one dev mode.")
